$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its textual representation (it contains
# locale-formatted numbers like "42.826.90" that Excel would otherwise
# coerce into numeric values, dropping the thousands separators).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.826.90"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").Value = "2.304.50"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "305.65"
$ws.Range("E5").Value = "  +2.19%  "

$ws.Range("D6").Value = "96.84"
$ws.Range("E6").Value = "  -0.35%  "

$ws.Range("E7").Value = "  -1.21%  "

$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("D9").Value = "0.498"
$ws.Range("E9").Value = "  -1.37%  "

$ws.Range("D10").Value = "35.35"
$ws.Range("E10").Value = "  -1.40%  "

$ws.Range("D11").Value = "0.0784"
$ws.Range("E11").Value = "  -0.32%  "

$ws.Range("D12").Value = "18.61"
$ws.Range("E12").Value = "  +4.67%  "

$ws.Range("E13").Value = "  +2.11%  "

$ws.Range("D14").Value = "6.85"
$ws.Range("E14").Value = "  +1.61%  "

$ws.Range("D15").Value = "2.657.60"
$ws.Range("E15").Value = "  +0.42%  "

$ws.Range("D16").Value = "2.262.77"
$ws.Range("E16").Value = "  -1.42%  "

$ws.Range("D17").Value = "0.780"
$ws.Range("E17").Value = "  +0.54%  "

$ws.Range("D18").Value = "42.754.36"
$ws.Range("E18").Value = "  +0.05%  "

$ws.Range("D19").Value = "12.69"
$ws.Range("E19").Value = "  +1.06%  "

$ws.Range("D20").Value = "0.0₃0894"
$ws.Range("E20").Value = "  -1.07%  "

$ws.Range("D21").Value = "6.02"
$ws.Range("E21").Value = "  -0.34%  "

$ws.Range("D22").Value = "67.22"
$ws.Range("E22").Value = "  -0.81%  "

$ws.Range("D23").Value = "235.93"
$ws.Range("E23").Value = "  -2.21%  "

$ws.Range("E24").Value = "  +0.66%  "

$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("D26").Value = "2.42"
$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").Value = "24.72"
$ws.Range("E27").Value = "  -1.64%  "

$ws.Range("D28").Value = "166.00"
$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("D29").Value = "2.05"
$ws.Range("E29").Value = "  +1.10%  "

$ws.Range("D30").Value = "9.04"
$ws.Range("E30").Value = "  +0.18%  "

$ws.Range("D31").Value = "33.04"
$ws.Range("E31").Value = "  +0.73%  "

$ws.Range("E32").Value = "  +0.09%  "

$ws.Range("D33").Value = "18.02"
$ws.Range("E33").Value = "  +5.37%  "

$ws.Range("E34").Value = "  -0.60%  "

$ws.Range("E35").Value = "  -6.10%  "

$ws.Range("E36").Value = "  -0.66%  "

$ws.Range("E37").Value = "  +0.27%  "

$ws.Range("E38").Value = "  +0.36%  "

$ws.Range("E39").Value = "  -0.36%  "

$ws.Range("E40").Value = "  -0.35%  "

$ws.Range("D41").Value = "2.71"
$ws.Range("E41").Value = "  -1.44%  "

$ws.Range("D42").Value = "2.000.10"
$ws.Range("E42").Value = "  -0.79%  "

$ws.Range("D43").Value = "0.0281"
$ws.Range("E43").Value = "  -0.91%  "

$ws.Range("D44").Value = "10.21"
$ws.Range("E44").Value = "  +1.15%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "18.04"
$ws.Range("E45").Value = "  +5.74%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "2.11"
$ws.Range("E46").Value = "  +1.08%  "

$ws.Range("D47").Value = "2.76"
$ws.Range("E47").Value = "  -0.30%  "

$ws.Range("D48").Value = "53.57"
$ws.Range("E48").Value = "  +0.91%  "

$ws.Range("D49").Value = "2.525.50"

$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").Value = "2.83"
$ws.Range("E50").Value = "  +1.61%  "

$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "71.07"
$ws.Range("E51").Value = "  -1.22%  "
